# Apply updated cryptocurrency market data (prices + 1h volume change)
# to the "cryptos" worksheet, matching the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "40.634.74"
$ws.Cells.Item(2, 5).Value = "  +2.74%  "
$ws.Cells.Item(3, 4).Value = "2.212.32"
$ws.Cells.Item(3, 5).Value = "  +1.68%  "
$ws.Cells.Item(4, 5).Value = "  +0.03%  "
$ws.Cells.Item(5, 4).Value = "'229.25"
$ws.Cells.Item(5, 5).Value = "  -0.31%  "
$ws.Cells.Item(6, 4).Value = "'0.632"
$ws.Cells.Item(6, 5).Value = "  +1.20%  "
$ws.Cells.Item(7, 4).Value = "'64.28"
$ws.Cells.Item(7, 5).Value = "  -1.34%  "
$ws.Cells.Item(9, 4).Value = "'0.406"
$ws.Cells.Item(9, 5).Value = "  +1.08%  "
$ws.Cells.Item(10, 4).Value = "'0.0865"
$ws.Cells.Item(10, 5).Value = "  -0.36%  "
$ws.Cells.Item(11, 5).Value = "  +0.27%  "
$ws.Cells.Item(12, 4).Value = "2.540.36"
$ws.Cells.Item(12, 5).Value = "  +1.66%  "
$ws.Cells.Item(13, 4).Value = "'15.94"
$ws.Cells.Item(13, 5).Value = "  -1.12%  "
$ws.Cells.Item(14, 4).Value = "'22.24"
$ws.Cells.Item(14, 5).Value = "  -1.16%  "
$ws.Cells.Item(15, 4).Value = "'0.826"
$ws.Cells.Item(15, 5).Value = "  +0.60%  "
$ws.Cells.Item(16, 5).Value = "  -0.01%  "
$ws.Cells.Item(17, 4).Value = "2.211.35"
$ws.Cells.Item(17, 5).Value = "  +2.58%  "
$ws.Cells.Item(18, 4).Value = "40.550.52"
$ws.Cells.Item(18, 5).Value = "  +2.65%  "
$ws.Cells.Item(19, 4).Value = "'73.82"
$ws.Cells.Item(19, 5).Value = "  +1.71%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0907"
$ws.Cells.Item(20, 5).Value = "  +5.59%  "
$ws.Cells.Item(21, 4).Value = "'6.12"
$ws.Cells.Item(21, 5).Value = "  -1.10%  "
$ws.Cells.Item(22, 4).Value = "'249.16"
$ws.Cells.Item(22, 5).Value = "  +7.02%  "
$ws.Cells.Item(23, 5).Value = "  +0.03%  "
$ws.Cells.Item(24, 4).Value = "'2.37"
$ws.Cells.Item(24, 5).Value = "  -0.39%  "
$ws.Cells.Item(25, 5).Value = "  -2.61%  "
$ws.Cells.Item(26, 4).Value = "'9.70"
$ws.Cells.Item(26, 5).Value = "  -0.58%  "
$ws.Cells.Item(27, 4).Value = "'173.35"
$ws.Cells.Item(27, 5).Value = "  +0.41%  "
$ws.Cells.Item(28, 4).Value = "'0.142"
$ws.Cells.Item(28, 5).Value = "  +1.58%  "
$ws.Cells.Item(29, 4).Value = "'20.40"
$ws.Cells.Item(29, 5).Value = "  +1.32%  "
$ws.Cells.Item(30, 4).Value = "'1.44"
$ws.Cells.Item(30, 5).Value = "  +2.28%  "
$ws.Cells.Item(31, 4).Value = "'2.82"
$ws.Cells.Item(31, 5).Value = "  +1.93%  "
$ws.Cells.Item(32, 4).Value = "'0.124"
$ws.Cells.Item(32, 5).Value = "  +0.83%  "
$ws.Cells.Item(33, 5).Value = "  -0.38%  "
$ws.Cells.Item(34, 4).Value = "'4.75"
$ws.Cells.Item(34, 5).Value = "  -1.85%  "
$ws.Cells.Item(35, 4).Value = "'7.07"
$ws.Cells.Item(35, 5).Value = "  -1.81%  "
# Row 36/37 swapped: Hedera now ranks above RenderToken
$ws.Cells.Item(36, 2).Value = "Hedera"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(36, 4).Value = "'0.0631"
$ws.Cells.Item(36, 5).Value = "  +1.15%  "
$ws.Cells.Item(37, 2).Value = "RenderToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(37, 4).Value = "'3.83"
$ws.Cells.Item(37, 5).Value = "  +6.10%  "
$ws.Cells.Item(38, 4).Value = "'2.47"
$ws.Cells.Item(38, 5).Value = "  +0.74%  "
$ws.Cells.Item(39, 4).Value = "'1.00"
$ws.Cells.Item(39, 5).Value = "  -0.05%  "
$ws.Cells.Item(40, 4).Value = "'4.78"
$ws.Cells.Item(40, 5).Value = "  +10.95%  "
$ws.Cells.Item(41, 4).Value = "'0.0233"
$ws.Cells.Item(41, 5).Value = "  +0.71%  "
$ws.Cells.Item(42, 4).Value = "'103.19"
$ws.Cells.Item(42, 5).Value = "  -1.84%  "
$ws.Cells.Item(43, 5).Value = "  +7.63%  "
# Rows 44-46 reordered: TrustWalletToken, InjectiveProtocol, Maker
$ws.Cells.Item(44, 2).Value = "TrustWalletToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(44, 4).Value = "'1.23"
$ws.Cells.Item(44, 5).Value = "  +2.91%  "
$ws.Cells.Item(45, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(45, 4).Value = "'17.54"
$ws.Cells.Item(45, 5).Value = "  -2.83%  "
$ws.Cells.Item(46, 2).Value = "Maker"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(46, 4).Value = "1.524.02"
$ws.Cells.Item(46, 5).Value = "  -1.12%  "
$ws.Cells.Item(47, 4).Value = "'1.11"
$ws.Cells.Item(47, 5).Value = "  -0.35%  "
$ws.Cells.Item(48, 4).Value = "'0.0933"
$ws.Cells.Item(48, 5).Value = "  +0.38%  "
$ws.Cells.Item(49, 4).Value = "'2.81"
$ws.Cells.Item(50, 4).Value = "'0.000205"
$ws.Cells.Item(50, 5).Value = "  +39.26%  "
$ws.Cells.Item(51, 4).Value = "'51.47"
$ws.Cells.Item(51, 5).Value = "  +9.84%  "
